# Auto-update draw results: append the 2025-10-19 Pick 3 draw as a new
# row at the bottom of the Results sheet (row 33), mirroring the existing
# rows 2-32 which hold plain text values in columns A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 33

# Columns A (date) and C (phase code) contain strings that look like a
# number/date ("2025-10-19", "251019"). Force the cells to the Text
# number format *before* assigning the value so they are stored as text
# (matching every other row in the sheet) instead of being coerced into
# a date serial / numeric value.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2025-10-19"

$ws.Range("B$newRow").Value = "Pick 3"

$ws.Range("C$newRow").NumberFormat = "@"
$ws.Range("C$newRow").Value = "251019"

$ws.Range("D$newRow").Value = "1-4-7"

$ws.Range("E$newRow").Value = "2025-10-19T21:36:16.785+04:00"
